$wb = $excel.ActiveWorkbook

$wsIn = $wb.Worksheets.Item("in")
$wsOut = $wb.Worksheets.Item("out")

# --- out!D2: edit the existing shared string text in place (服务器首月 -> 服务器首周) ---
$wsOut.Cells.Item(2, 4).Value = "服务器首周"

# --- in sheet: C24 amount correction ---
$wsIn.Cells.Item(24, 3).Value = 61

# --- out sheet row 9: shift the previous row 8 transaction down to row 9 ---
$wsOut.Cells.Item(9, 1).Value = 20190912
$wsOut.Cells.Item(9, 2).Value = "ZL Asica"
$wsOut.Cells.Item(9, 3).Value = 80.099999999999994
$wsOut.Cells.Item(9, 4).Value = "服务器9月费用"

# --- out sheet row 8: new transaction (first new shared string: index 47) ---
$wsOut.Cells.Item(8, 1).Value = 20190901
$wsOut.Cells.Item(8, 2).Value = "keishi"
$wsOut.Cells.Item(8, 3).Value = 41
$wsOut.Cells.Item(8, 4).Value = "半月服务器费用"

# --- in sheet rows 26-27: new transactions (shared strings 48, 49) ---
$wsIn.Cells.Item(26, 1).Value = 20191005
$wsIn.Cells.Item(26, 2).Value = "FAKED"
$wsIn.Cells.Item(26, 3).Value = 5
$wsIn.Cells.Item(26, 4).Value = "wechat"

$wsIn.Cells.Item(27, 1).Value = 20191007
$wsIn.Cells.Item(27, 2).Value = "*尔"
$wsIn.Cells.Item(27, 3).Value = 100
$wsIn.Cells.Item(27, 4).Value = "wechat"

# --- out sheet row 10: new transaction (shared strings 50, 51) ---
$wsOut.Cells.Item(10, 1).Value = 20191003
$wsOut.Cells.Item(10, 2).Value = "Keishi"
$wsOut.Cells.Item(10, 3).Value = 98
$wsOut.Cells.Item(10, 4).Value = "B站年度大会员"

# --- in sheet row 28: new transaction (shared string 52) ---
$wsIn.Cells.Item(28, 1).Value = 20191007
$wsIn.Cells.Item(28, 2).Value = "*翔辰"
$wsIn.Cells.Item(28, 3).Value = 3
$wsIn.Cells.Item(28, 4).Value = "alipay"

# --- out sheet: append a new blank styled row 24 (copy formatting from row 23) ---
$wsOut.Range("A23:C23").Copy($wsOut.Range("A24:C24"))

# --- view state ---
$wsOut.Activate()
$wsOut.Range("E13").Select()

$wsIn.Activate()
$wsIn.Range("D28").Select()
